$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("erosion")
$ws2 = $wb.Worksheets.Item("soil")

# ---------------------------------------------------------------------------
# Sheet "erosion": replace the Beijing erosion-rate series with the Dunhuang
# series, add a new 10th data point (row 11), and resize the label columns.
# ---------------------------------------------------------------------------
$erosionData = @(
    @(0, 0.32500000000000001, 0.1,  "敦煌 "),
    @(1, 0.314,                1.5,  "敦煌 "),
    @(2, 0.42899999999999999, 3,    "敦煌 "),
    @(3, 0.40100000000000002, 4.5,  "敦煌 "),
    @(4, 0.45900000000000002, 6,    "敦煌 "),
    @(5, 0.435,                7.5,  "敦煌 "),
    @(6, 0.53,                 9,    "敦煌 "),
    @(7, 0.499,                10.5, "敦煌 "),
    @(8, 0.44900000000000001, 12,   "敦煌 "),
    @(9, 0.59899999999999998, 18,   "敦煌 ")
)

for ($i = 0; $i -lt $erosionData.Count; $i++) {
    $row = 2 + $i
    $vals = $erosionData[$i]
    $ws1.Range("A$row").Value = $vals[0]
    $ws1.Range("B$row").Value = $vals[1]
    $ws1.Range("C$row").Value = $vals[2]
    $ws1.Range("D$row").Value = $vals[3]
}

$ws1.Columns.Item(2).ColumnWidth = 32
$ws1.Columns.Item(3).ColumnWidth = 29.714285714285715

# ---------------------------------------------------------------------------
# Sheet "soil": replace the single Beijing climate-summary row with the
# Dunhuang values.
# ---------------------------------------------------------------------------
$ws2.Range("B2").Value = 11.3
$ws2.Range("C2").Value = 38
$ws2.Range("D2").Value = 48
$ws2.Range("E2").Value = 7
$ws2.Range("F2").Value = 2632
$ws2.Range("G2").Value = 7846
$ws2.Range("H2").Value = "敦煌 "

$ws2.Columns.Item(7).ColumnWidth = 8.142857142857142

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to mirror the saved view state.
# ---------------------------------------------------------------------------
$ws2.Range("G12").Select()
$ws1.Activate()
$ws1.Range("A11").Select()

Write-Output "done"
